$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 3
$ws.Range("R2").Value = 1.7
$ws.Range("T2").Value = 1.25
$ws.Range("AO2").Value = 34
